$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 7577703.5
$ws.Range("I106").Value = 11905850
$ws.Range("J106").Value = 3447.25
$ws.Range("K106").Value = 11905850
$ws.Range("L106").Value = 3447.25
$ws.Range("M106").Value = -11905219
$ws.Range("N106").Value = -4709.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 555.1667
$ws.Range("I125").Value = 351.7143
$ws.Range("K125").Value = 3165.4287
$ws.Range("M125").Value = -705.4286999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2553.7222
$ws.Range("I132").Value = 2754.25
$ws.Range("J132").Value = 949.5
$ws.Range("K132").Value = 8262.75
$ws.Range("L132").Value = 2848.5
$ws.Range("M132").Value = -5732.75
$ws.Range("N132").Value = -7908.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 57277.945
$ws.Range("I137").Value = 1614.5714
$ws.Range("J137").Value = 92700.09
$ws.Range("K137").Value = 4843.7142
$ws.Range("L137").Value = 278100.27
$ws.Range("M137").Value = -2293.7142
$ws.Range("N137").Value = -283200.27

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1220.0731
$ws.Range("I141").Value = 883.5833
$ws.Range("J141").Value = 3642.8
$ws.Range("K141").Value = 2650.7499
$ws.Range("L141").Value = 10928.4
$ws.Range("M141").Value = 2529.2501
$ws.Range("N141").Value = -21288.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 20332.666
$ws.Range("J24").Value = 20332.666
$ws.Range("L24").Value = 20332.666
$ws.Range("N24").Value = -21080.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2244.8
$ws.Range("J45").Value = 2587.2727
$ws.Range("L45").Value = 2587.2727
$ws.Range("N45").Value = -3341.2727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 14296.777
$ws.Range("J92").Value = 14296.777
$ws.Range("L92").Value = 14296.777
$ws.Range("N92").Value = -19288.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 20332.666
$ws.Range("J100").Value = 20332.666
$ws.Range("L100").Value = 20332.666
$ws.Range("N100").Value = -22496.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1270.8572
$ws.Range("I20").Value = 1340.5
$ws.Range("K20").Value = 1340.5
$ws.Range("M20").Value = -1093.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 27662.293
$ws.Range("I134").Value = 29530.37
$ws.Range("K134").Value = 88591.11
$ws.Range("M134").Value = -86056.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10615.359
$ws.Range("I31").Value = 10984.344
$ws.Range("J31").Value = 8928.571
$ws.Range("K31").Value = 10984.344
$ws.Range("L31").Value = 8928.571
$ws.Range("M31").Value = -10689.344
$ws.Range("N31").Value = -9518.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10615.359
$ws.Range("I34").Value = 10984.344
$ws.Range("J34").Value = 8928.571
$ws.Range("K34").Value = 10984.344
$ws.Range("L34").Value = 8928.571
$ws.Range("M34").Value = -10782.344
$ws.Range("N34").Value = -9332.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 21419.56
$ws.Range("I58").Value = 1442.7222
$ws.Range("J58").Value = 72788.57000000001
$ws.Range("K58").Value = 1442.7222
$ws.Range("L58").Value = 72788.57000000001
$ws.Range("M58").Value = -1239.7222
$ws.Range("N58").Value = -73194.57000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 19332.666
$ws.Range("J92").Value = 19332.666
$ws.Range("L92").Value = 19332.666
$ws.Range("N92").Value = -24324.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 30666.666
$ws.Range("J95").Value = 30666.666
$ws.Range("L95").Value = 30666.666
$ws.Range("N95").Value = -36158.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1546.7826
$ws.Range("I107").Value = 1251.8334
$ws.Range("J107").Value = 1650.8823
$ws.Range("K107").Value = 1251.8334
$ws.Range("L107").Value = 1650.8823
$ws.Range("M107").Value = 668.1666
$ws.Range("N107").Value = -5490.8823

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 16017.473
$ws.Range("I132").Value = 15874.714
$ws.Range("K132").Value = 47624.142
$ws.Range("M132").Value = -45094.142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 730.25
$ws.Range("I134").Value = 742.7143
$ws.Range("K134").Value = 2228.1429
$ws.Range("M134").Value = 306.8571000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 21419.56
$ws.Range("I136").Value = 1442.7222
$ws.Range("J136").Value = 72788.57000000001
$ws.Range("K136").Value = 4328.1666
$ws.Range("L136").Value = 218365.71
$ws.Range("M136").Value = -1778.1666
$ws.Range("N136").Value = -223465.71

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11645.958
$ws.Range("I70").Value = 12951.6
$ws.Range("J70").Value = 10713.357
$ws.Range("K70").Value = 12951.6
$ws.Range("L70").Value = 10713.357
$ws.Range("M70").Value = -12681.6
$ws.Range("N70").Value = -11253.357

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 11645.958
$ws.Range("I73").Value = 12951.6
$ws.Range("J73").Value = 10713.357
$ws.Range("K73").Value = 12951.6
$ws.Range("L73").Value = 10713.357
$ws.Range("M73").Value = -12015.6
$ws.Range("N73").Value = -12585.357

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3587.1304
$ws.Range("I80").Value = 2854.3635
$ws.Range("J80").Value = 4258.8335
$ws.Range("K80").Value = 2854.3635
$ws.Range("L80").Value = 4258.8335
$ws.Range("M80").Value = -1856.3635
$ws.Range("N80").Value = -6254.8335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3587.1304
$ws.Range("I83").Value = 2854.3635
$ws.Range("J83").Value = 4258.8335
$ws.Range("K83").Value = 14271.8175
$ws.Range("L83").Value = 21294.1675
$ws.Range("M83").Value = -9279.817499999999
$ws.Range("N83").Value = -31278.1675

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4528.8125
$ws.Range("I126").Value = 3714.6667
$ws.Range("K126").Value = 11144.0001
$ws.Range("M126").Value = -8674.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 47624270
$ws.Range("I7").Value = 71431450
$ws.Range("J7").Value = 9915.429
$ws.Range("K7").Value = 71431450
$ws.Range("L7").Value = 9915.429
$ws.Range("M7").Value = -71431338
$ws.Range("N7").Value = -10139.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1470.6
$ws.Range("I93").Value = 1356.2222
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1356.2222
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -108.2221999999999
$ws.Range("N93").Value = -4996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 8339.666999999999
$ws.Range("J94").Value = 8339.666999999999
$ws.Range("L94").Value = 8339.666999999999
$ws.Range("N94").Value = -9691.666999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 786903.75
$ws.Range("I122").Value = 1510407.2
$ws.Range("K122").Value = 4531221.6
$ws.Range("M122").Value = -4528771.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 47624270
$ws.Range("I126").Value = 71431450
$ws.Range("J126").Value = 9915.429
$ws.Range("K126").Value = 214294350
$ws.Range("L126").Value = 29746.287
$ws.Range("M126").Value = -214291880
$ws.Range("N126").Value = -34686.287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 15734.912
$ws.Range("I136").Value = 17562.4
$ws.Range("K136").Value = 52687.2
$ws.Range("M136").Value = -50137.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4385.857
$ws.Range("I62").Value = 4001
$ws.Range("J62").Value = 4539.8
$ws.Range("K62").Value = 4001
$ws.Range("L62").Value = 4539.8
$ws.Range("M62").Value = -3377
$ws.Range("N62").Value = -5787.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4385.857
$ws.Range("I65").Value = 4001
$ws.Range("J65").Value = 4539.8
$ws.Range("K65").Value = 20005
$ws.Range("L65").Value = 22699
$ws.Range("M65").Value = -16885
$ws.Range("N65").Value = -28939

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 18091.75
$ws.Range("J104").Value = 18091.75
$ws.Range("L104").Value = 18091.75
$ws.Range("N104").Value = -25079.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1161.1666
$ws.Range("I126").Value = 1124.9
$ws.Range("J126").Value = 1342.5
$ws.Range("K126").Value = 3374.7
$ws.Range("L126").Value = 4027.5
$ws.Range("M126").Value = -904.7000000000003
$ws.Range("N126").Value = -8967.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 580.9286
$ws.Range("I132").Value = 543.5641000000001
$ws.Range("J132").Value = 1066.6666
$ws.Range("K132").Value = 1630.6923
$ws.Range("L132").Value = 3199.9998
$ws.Range("M132").Value = 899.3076999999998
$ws.Range("N132").Value = -8259.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 22941514
$ws.Range("I136").Value = 27165726
$ws.Range("J136").Value = 10085.714
$ws.Range("K136").Value = 81497178
$ws.Range("L136").Value = 30257.142
$ws.Range("M136").Value = -81494628
$ws.Range("N136").Value = -35357.142
